$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.873.78"
$ws.Range("E2").Value = "  +1.52%  "

# Row 3
$ws.Range("D3").Value = "2.298.31"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.71"
$ws.Range("E5").Value = "  -1.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.98"
$ws.Range("E6").Value = "  +1.77%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -0.41%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.24"
$ws.Range("E10").Value = "  -0.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -0.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.30"
$ws.Range("E12").Value = "  -0.01%  "

# Row 13
$ws.Range("E13").Value = "  +1.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.992"
$ws.Range("E14").Value = "  +2.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.16"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("D16").Value = "2.648.27"
$ws.Range("E16").Value = "  +0.26%  "

# Row 17
$ws.Range("D17").Value = "2.294.08"
$ws.Range("E17").Value = "  -0.20%  "

# Row 18
$ws.Range("D18").Value = "42.490.28"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  -1.40%  "

# Row 20
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.68"
$ws.Range("E20").Value = "  +5.09%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("E21").Value = "  -0.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.32"
$ws.Range("E22").Value = "  +0.52%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.52"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.60"
$ws.Range("E24").Value = "  -1.77%  "

# Row 25
$ws.Range("E25").Value = "  -2.00%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.76"
$ws.Range("E27").Value = "  -0.89%  "

# Row 28
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.00"
$ws.Range("E28").Value = "  +14.22%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("E29").Value = "  -4.04%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.42"
$ws.Range("E30").Value = "  -1.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.10"
$ws.Range("E31").Value = "  -3.99%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.29"
$ws.Range("E32").Value = "  -0.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0865"
$ws.Range("E33").Value = "  -1.17%  "

# Row 34
$ws.Range("E34").Value = "  -1.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.58"
$ws.Range("E35").Value = "  -0.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("E37").Value = "  -1.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0350"
$ws.Range("E38").Value = "  -1.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.76"
$ws.Range("E39").Value = "  +1.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("E40").Value = "  -1.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.60"
$ws.Range("E41").Value = "  +3.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.07"
$ws.Range("E42").Value = "  +5.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.44"
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.227"
$ws.Range("E44").Value = "  +1.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.43%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.22"
$ws.Range("E46").Value = "  +1.75%  "

# Row 47
$ws.Range("D47").Value = "1.745.16"
$ws.Range("E47").Value = "  +8.52%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "80.30"
$ws.Range("E48").Value = "  +1.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "111.68"
$ws.Range("E49").Value = "  -3.74%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.19"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.70"
$ws.Range("E51").Value = "  -3.02%  "
